# Generate Report for Handback
# Update the "Latest Handback DateTime" and "Error Detail" values for the
# 614fb894-1524-4b63-a580-394b1a45dc50 row on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newErrorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a54b6e461f7f27fa510be098785493daea29/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40d86b4d7842a7cedbf32bbd9c64b7886a1c84ce/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md."

# zh-cn sheet, row 3 (614fb894... file)
$wsZhCn.Range("L3").Value = "2017-02-17 09:39:25"
$wsZhCn.Range("R3").Value = $newErrorDetail

# de-de sheet, row 3 (614fb894... file)
$wsDeDe.Range("L3").Value = "2017-02-17 09:39:48"
$wsDeDe.Range("R3").Value = $newErrorDetail
